# email sent on 06022020
# Duplicate rows 2-8 of Sheet1 into rows 9-15 (same people/emails, new
# usernames + attendance count), and append the two matching summary rows
# to Sheet2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Source rows (2-8) mapped onto destination rows (9-15), carrying the same
# Email (column D) value as the original row, but a brand-new random
# Username (column J) and a bumped Attendance (column K) value.
$rows = @(
    @{ Src = 2; Dst = 9;  Email = "sourabh.awasthi@capgemini.com";  User = "&Np=+4(J2}3z" },
    @{ Src = 3; Dst = 10; Email = "sandipan.deb@capgemini.com";     User = "w}xY%3i(gVU8" },
    @{ Src = 4; Dst = 11; Email = "biswaji.deb@capgemini.com";      User = "lUI%Uxy@twfW" },
    @{ Src = 5; Dst = 12; Email = "debanjan.das@capgemini.com";     User = "ZC4/UCH+H//x" },
    @{ Src = 6; Dst = 13; Email = "dhiraj.kajari@capgemini.com";    User = "3xkYs}F7UORG" },
    @{ Src = 7; Dst = 14; Email = "manoj-kumar.b.s@capgemini.com";  User = "a}DB@fw]ZBfX" },
    @{ Src = 8; Dst = 15; Email = "mayur.bhorkar@capgemini.com";    User = "9(fZ)e6AV6M]" }
)

$q = [char]34

foreach ($row in $rows) {
    $s = $row.Src
    $d = $row.Dst

    $ws1.Range("A$d").Formula = "=PROPER(IFERROR(LEFT(C$d,FIND(CHAR(46),C$d)-1),C$d))"
    $ws1.Range("B$d").Formula = "=IFERROR(PROPER(RIGHT(C$d,LEN(C$d)-FIND(${q}@${q},SUBSTITUTE(C$d,${q}.${q},${q}@${q},((LEN(C$d)-LEN(SUBSTITUTE(C$d,${q}.${q},${q}${q})))/LEN(${q}\${q})))))), ${q}Unknown${q})"
    $ws1.Range("C$d").Formula = "=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D$d,FIND(CHAR(64),D$d)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))"
    $ws1.Range("D$d").Value = $row.Email
    $ws1.Range("E$d").Formula = "=LEFT(H$d,FIND(CHAR(46),H$d)-1)"
    $ws1.Range("F$d").Formula = "=CONCATENATE(${q}ITPartner\${q},I$d)"
    $ws1.Range("H$d").Formula = "=RIGHT(D$d,LEN(D$d)-FIND(CHAR(64),D$d))"
    $ws1.Range("I$d").Formula = "=PROPER(E$d)"
    $ws1.Range("J$d").Value = $row.User
    $ws1.Range("K$d").Value = 81
    $ws1.Range("M$d").Value = $true
    $ws1.Range("P$d").Formula = "=COUNTIF(D:D,D$d)"
}

# Sheet2 gains two more summary rows echoing the (already-duplicated)
# sandipan.deb and mayur.bhorkar entries.
$ws2.Range("C4").Value = "sandipan.deb"
$ws2.Range("I4").Value = "Capgemini"
$ws2.Range("C5").Value = "mayur.bhorkar"
$ws2.Range("I5").Value = "Capgemini"
